# Insert a new row at row 222, shifting existing rows 222:298 down to 223:299.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with the new record's data.
# Columns that stay identical to the template (copied down from the row that
# used to be at 222, now at 223) are left as Excel set them via the insert
# (A,B,C,E,F,G,H,I,J,K,L,Q,R,T); only the record-specific fields are set here.
$ws.Range("A222").Value = 10
$ws.Range("B222").Value = "Vega Modelo de Temuco"
$ws.Range("C222").Value = "La Araucanía"
$ws.Range("D222").Value = 44809
$ws.Range("E222").Value = 9
$ws.Range("F222").Value = "Fruta"
$ws.Range("G222").Value = 100102
$ws.Range("H222").Value = "Cítricos"
$ws.Range("I222").Value = 100102006
$ws.Range("J222").Value = "Pomelo"
$ws.Range("K222").Value = "Start Ruby"
$ws.Range("L222").Value = "Primera"
$ws.Range("M222").Value = 125
$ws.Range("N222").Value = 13000
$ws.Range("O222").Value = 13000
$ws.Range("P222").Value = 13000
$ws.Range("Q222").Value = "$/bandeja 15 kilos granel"
$ws.Range("R222").Value = "Región de O'Higgins"
$ws.Range("S222").Value = 867
$ws.Range("T222").Value = 15
